$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Zebronics - Wireless Bluetooth Headset"
$ws.Range("B1").Value = "Rs. 1,449"

$ws.Range("A2").Value = "BLUNT Air 1 Earbuds In Ear Bluetooth Earphone 6 Hours Playback Bluetooth IPX6(Water Resistant) Active Noise cancellation -Bluetooth V 5.1 Black"
$ws.Range("B2").Value = "Rs. 1,299"

$ws.Range("A3").Value = "boAt Airdopes 131/138 Twin Wireless Earbuds with IWP Technology, Bluetooth V5.0, Immersive Audio, Up to 15H Total Playback, Instant Voice Assistant and Type-C Charging,Bluetooth Earphone (Active Black)"
$ws.Range("B3").Value = "Rs. 1,199"

$ws.Range("A4").Value = "hitage TWS68 V5.0Earbuds In Ear True Wireless (TWS) 20 Hours Playback IPX4(Splash & Sweat Proof) Comfirtable in ear fit -Bluetooth V 5.0 Red"
$ws.Range("B4").Value = "Rs. 795"

$ws.Range("A5").Value = "Tecsox PowerHouse Earbud In Ear Bluetooth Earphone 45 Hours Playback Bluetooth IPX5(Splash Proof) Powerfull Bass -Bluetooth V 5.1 Black"
$ws.Range("B5").Value = "Rs. 725"
